$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (changed) date in column C for all data rows (2-45)
# from 2023-09-15 (45184) to 2023-09-17 (45186).
for ($r = 2; $r -le 45; $r++) {
    $ws.Cells.Item($r, 3).Value = 45186
}

# Add a friendly display-text second argument ("Beteckning") to the
# HYPERLINK formulas in columns S, T, V, W, X, Y for rows 2-6.
$segments = @{
    "S" = "artfynd"
    "T" = "kartor"
    "V" = "klagomål"
    "W" = "klagomålsmail"
    "X" = "tillsyn"
    "Y" = "tillsynsmail"
}
$exts = @{
    "S" = ".xlsx"
    "T" = ".png"
    "V" = ".docx"
    "W" = ".docx"
    "X" = ".docx"
    "Y" = ".docx"
}

for ($r = 2; $r -le 6; $r++) {
    $beteckning = $ws.Range("A$r").Value2
    foreach ($col in @("S", "T", "V", "W", "X", "Y")) {
        $url = "https://klasma.github.io/Logging_TJORN/" + $segments[$col] + "/" + $beteckning + $exts[$col]
        $formula = '=HYPERLINK("' + $url + '", "' + $beteckning + '")'
        $ws.Range("$col$r").Formula = $formula
    }
}
